$p = $ppt.ActivePresentation
$nm = $p.NotesMaster
Write-Host "NotesMaster props:"
Write-Host "  Height:" $nm.Height
Write-Host "  Width:" $nm.Width
$bg = $nm.Background
Write-Host "  Background:" $bg
